# edit.ps1
#
# Commit: "Added Saving/Loading. Moved FPS control to TimeController.
# Added ProcessManager. Switched from window.setInterval(this.loop) to
# requestAnimationFrame(this.loop). Added Settings. Changed combat canvas
# size from 1600x1600 to 960x960. Added createEnum() and enums. Created
# TextureLoader."
#
# The visible document change is a new "ORDER:" planning/checklist block
# appended to the very end of the design doc, after the last (blank)
# paragraph and before the section properties.
#
# We build the whole block as one WordprocessingML fragment and insert it
# with Range.InsertXML so that:
#   * the blank separator paragraphs come out as clean self-closing
#     <w:p/> elements (Range.InsertParagraphAfter / a typed Enter leaves a
#     stray empty <w:r/> behind instead, which would not match);
#   * the "option" paragraph keeps its three separate <w:r> runs
#     ("Block health makes the block fade and opt" + "i" + "on").
#
# Quirk: InsertXML at a collapsed Range sitting exactly on an existing
# paragraph boundary treats the FIRST <w:p/> in the fragment as that
# boundary (no new paragraph is created for it) -- only the <w:p/>
# elements after it produce genuinely new paragraphs. Since the target
# adds two brand-new blank paragraphs, the fragment below starts with
# three <w:p/> elements.

$d = $word.ActiveDocument

# Last paragraph in the body -- the trailing empty <w:p/> right before sectPr.
$lastParaIndex = $d.Paragraphs.Count
$rng = $d.Paragraphs($lastParaIndex).Range
$rng.Collapse(0)

$fragment = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>ORDER:</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Settings</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Save settings/settings profile</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Block health display settings</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Block health display real angle vs fast angle</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Block health display health left to right</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Option to not display any of them</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Block health percentage and option</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Can have any or none of the health displays</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Block health makes the block fade and opt</w:t></w:r><w:r><w:t>i</w:t></w:r><w:r><w:t>on</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>UIState and UIPanel ShowHideUI</w:t></w:r></w:p>'

[void]$rng.InsertXML($fragment)
